$wb = $excel.ActiveWorkbook

# Sheet "展览": F2 6874 -> 6876, F5 35 -> 36
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 6876
$wsExhibit.Range("F5").Value = 36

# Sheet "全部类型": F2 6874 -> 6876, F5 35 -> 36
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 6876
$wsAll.Range("F5").Value = 36
